$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate the header row (row 1) from English to Russian ---
$ws.Range("B1").Value = "имя"
$ws.Range("C1").Value = "фамилия"
$ws.Range("D1").Value = "отчество"
$ws.Range("E1").Value = "возраст"
$ws.Range("F1").Value = "телефон"
$ws.Range("I1").Value = "улица"
$ws.Range("J1").Value = "почта"
$ws.Range("K1").Value = "пол"
$ws.Range("L1").Value = "паспорт"

# --- Add a new "должность" (job title) column in M ---
$ws.Range("M1").Value = "должность"

$titles = @("Стажёр", "Инженер", "Системный администратор", "Бухгалтер", "Паркхмахер", "Диспетчер", "Оленевод", "Охраник", "Директор")
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $titles[$i]
}

# --- Move the active selection to B1 ---
$ws.Range("B1").Select() | Out-Null
